$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell T2's value: 204387 -> 204671
$ws.Range("T2").Value = 204671

# Move the active selection from T2 to T3 (no data change, just cursor move)
$ws.Range("T3").Select()
